# Remove the full worked-solution body of exercise 09-34 (parts (1)-(4)
# plus the closing remark), leaving only the "Solution:" heading and the
# trailing blank paragraphs that originally followed the solution.

$d = $word.ActiveDocument
$paras = $d.Paragraphs

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $text = $p.Range.Text

    if ($null -eq $startPara -and $text -eq "(1)`r") {
        $startPara = $p
    }
    if ($text -like "The company is performing well*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
